$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
  # Row 17: One for the Road
$ws.Range("H17").Value = 6063.5625
$ws.Range("J17").Value = 6341.2
$ws.Range("L17").Value = 19023.6
$ws.Range("N17").Value = -19359.6
  # Row 62: The Mustache Suits Him
$ws.Range("H62").Value = 8111.375
$ws.Range("J62").Value = 8815.166999999999
$ws.Range("L62").Value = 8815.166999999999
$ws.Range("N62").Value = -10063.167
  # Row 65: Forgery of Convenience (L)
$ws.Range("H65").Value = 8111.375
$ws.Range("J65").Value = 8815.166999999999
$ws.Range("L65").Value = 44075.835
$ws.Range("N65").Value = -50315.835
  # Row 106: Making Your Mark
$ws.Range("H106").Value = 7248439
$ws.Range("I106").Value = 10102446
$ws.Range("K106").Value = 10102446
$ws.Range("M106").Value = -10101815
  # Row 107: Another Man's Ink
$ws.Range("H107").Value = 1134.4375
$ws.Range("I107").Value = 786.0714
$ws.Range("J107").Value = 3573
$ws.Range("K107").Value = 786.0714
$ws.Range("L107").Value = 3573
$ws.Range("M107").Value = 1133.9286
$ws.Range("N107").Value = -7413
  # Row 111: An Eye for Healing
$ws.Range("H111").Value = 3219.3333
$ws.Range("I111").Value = 4461.25
$ws.Range("J111").Value = 1800
$ws.Range("K111").Value = 13383.75
$ws.Range("L111").Value = 5400
$ws.Range("M111").Value = -10316.75
$ws.Range("N111").Value = -11534
  # Row 116: Growing Up
$ws.Range("H116").Value = 5012.875
$ws.Range("I116").Value = 2275
$ws.Range("K116").Value = 2275
$ws.Range("M116").Value = 1167
  # Row 118: Crafty Concoctions
$ws.Range("H118").Value = 750
$ws.Range("I118").Value = 750
$ws.Range("K118").Value = 2250
$ws.Range("M118").Value = -593
  # Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 4512.0435
$ws.Range("I132").Value = 4514.263
$ws.Range("J132").Value = 4501.5
$ws.Range("K132").Value = 13542.789
$ws.Range("L132").Value = 13504.5
$ws.Range("M132").Value = -11012.789
$ws.Range("N132").Value = -18564.5
  # Row 135: For Tired Minds
$ws.Range("H135").Value = 17243402
$ws.Range("J135").Value = 100009010
$ws.Range("L135").Value = 900081090
$ws.Range("N135").Value = -900086160
  # Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 1639.2
$ws.Range("I137").Value = 1457.375
$ws.Range("K137").Value = 4372.125
$ws.Range("M137").Value = -1822.125
  # Row 138: All-night Crafting
$ws.Range("H138").Value = 28574532
$ws.Range("I138").Value = 62501624
$ws.Range("J138").Value = 4349.316
$ws.Range("K138").Value = 187504872
$ws.Range("L138").Value = 13047.948
$ws.Range("M138").Value = -187499732
$ws.Range("N138").Value = -23327.948

$ws = $wb.Worksheets.Item("ARM")
  # Row 2: Ain't Got No Ingots
$ws.Range("H2").Value = 1528.3529
$ws.Range("I2").Value = 1362.1333
$ws.Range("K2").Value = 1362.1333
$ws.Range("M2").Value = -1249.1333
  # Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 3380.077
$ws.Range("I61").Value = 3235.28
$ws.Range("K61").Value = 3235.28
$ws.Range("M61").Value = -3023.28
  # Row 80: A Squire to Inspire
$ws.Range("H80").Value = 65535
$ws.Range("J80").Value = 65535
$ws.Range("L80").Value = 65535
$ws.Range("N80").Value = -67531
  # Row 83: All's Fair in Highborn Assassination (L)
$ws.Range("H83").Value = 65535
$ws.Range("J83").Value = 65535
$ws.Range("L83").Value = 196605
$ws.Range("N83").Value = -206589
  # Row 116: No Scope
$ws.Range("H116").Value = 1528.3529
$ws.Range("I116").Value = 1362.1333
$ws.Range("K116").Value = 1362.1333
$ws.Range("M116").Value = 931.8667
  # Row 126: Armoire Aftercare
$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 3000
$ws.Range("K126").Value = 9000
$ws.Range("M126").Value = -6530
  # Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 28349.3
$ws.Range("I132").Value = 3503
$ws.Range("K132").Value = 10509
$ws.Range("M132").Value = -7979
  # Row 136: Metal with Mettle
$ws.Range("H136").Value = 3380.077
$ws.Range("I136").Value = 3235.28
$ws.Range("K136").Value = 9705.84
$ws.Range("M136").Value = -7155.84

$ws = $wb.Worksheets.Item("BSM")
  # Row 3: Hells Bells
$ws.Range("H3").Value = 1528.3529
$ws.Range("I3").Value = 1362.1333
$ws.Range("K3").Value = 1362.1333
$ws.Range("M3").Value = -1248.1333
  # Row 105: Ingot to Wing It
$ws.Range("H105").Value = 3336705.2
$ws.Range("I105").Value = 3747.9
$ws.Range("K105").Value = 3747.9
$ws.Range("M105").Value = -2000.9

$ws = $wb.Worksheets.Item("CRP")
  # Row 4: A Clogful of Camaraderie
$ws.Range("H4").Value = 2890
$ws.Range("J4").Value = 2890
$ws.Range("L4").Value = 2890
$ws.Range("N4").Value = -3114
  # Row 7: Gridania's Got Talent
$ws.Range("H7").Value = 31.4
$ws.Range("I7").Value = 40
$ws.Range("J7").Value = 29.25
$ws.Range("K7").Value = 40
$ws.Range("L7").Value = 29.25
$ws.Range("M7").Value = 73
$ws.Range("N7").Value = -255.25
  # Row 31: Wall Not Found
$ws.Range("H31").Value = 2617.3572
$ws.Range("I31").Value = 1179.0769
$ws.Range("J31").Value = 3262.1035
$ws.Range("K31").Value = 1179.0769
$ws.Range("L31").Value = 3262.1035
$ws.Range("M31").Value = -884.0769
$ws.Range("N31").Value = -3852.1035
  # Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 2617.3572
$ws.Range("I34").Value = 1179.0769
$ws.Range("J34").Value = 3262.1035
$ws.Range("K34").Value = 1179.0769
$ws.Range("L34").Value = 3262.1035
$ws.Range("M34").Value = -977.0769
$ws.Range("N34").Value = -3666.1035
  # Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 1222.25
$ws.Range("I134").Value = 1076.9375
$ws.Range("K134").Value = 3230.8125
$ws.Range("M134").Value = -695.8125

$ws = $wb.Worksheets.Item("CUL")
  # Row 4: In Hot Water
$ws.Range("H4").Value = 941.25
$ws.Range("I4").Value = 941.25
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 2823.75
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -2711.75
$ws.Range("N4").Value = $null
  # Row 23: Sweet Smell of Success
$ws.Range("H23").Value = 441.75
$ws.Range("I23").Value = 39.666668
$ws.Range("K23").Value = 119.000004
$ws.Range("M23").Value = 115.999996
  # Row 80: Saucy for a Suitor
$ws.Range("H80").Value = 8427.143
$ws.Range("J80").Value = 8844.691999999999
$ws.Range("L80").Value = 26534.076
$ws.Range("N80").Value = -28406.076
  # Row 83: Saved by the Sauce (L)
$ws.Range("H83").Value = 8427.143
$ws.Range("J83").Value = 8844.691999999999
$ws.Range("L83").Value = 79602.22799999999
$ws.Range("N83").Value = -88962.22799999999
  # Row 113: Can't Eat Just One
$ws.Range("H113").Value = 437.5
$ws.Range("J113").Value = 443.33334
$ws.Range("L113").Value = 1330.00002
$ws.Range("N113").Value = -5670.000019999999
  # Row 131: The Mountain Steeped
$ws.Range("H131").Value = 700.63
$ws.Range("J131").Value = 718.2043
$ws.Range("L131").Value = 2154.6129
$ws.Range("N131").Value = -12234.6129
  # Row 136: Simple Is Hardest
$ws.Range("H136").Value = 2789.6365
$ws.Range("J136").Value = 4539.2
$ws.Range("L136").Value = 13617.6
$ws.Range("N136").Value = -23817.6

$ws = $wb.Worksheets.Item("GSM")
  # Row 5: Hora at Me
$ws.Range("H5").Value = 4933.1665
$ws.Range("I5").Value = 3399.75
$ws.Range("J5").Value = 8000
$ws.Range("K5").Value = 3399.75
$ws.Range("L5").Value = 8000
$ws.Range("M5").Value = -3287.75
$ws.Range("N5").Value = -8224
  # Row 70: Sky Is the Limit
$ws.Range("H70").Value = 3298261.2
$ws.Range("I70").Value = 3863.3333
$ws.Range("K70").Value = 3863.3333
$ws.Range("M70").Value = -3593.3333
  # Row 73: Hulls of Broken Dreams (L)
$ws.Range("H73").Value = 3298261.2
$ws.Range("I73").Value = 3863.3333
$ws.Range("K73").Value = 3863.3333
$ws.Range("M73").Value = -2927.3333
  # Row 107: Whetstones for the Workers
$ws.Range("H107").Value = 4273743
$ws.Range("I107").Value = 242.5
$ws.Range("J107").Value = 12820744
$ws.Range("K107").Value = 242.5
$ws.Range("L107").Value = 12820744
$ws.Range("M107").Value = 1677.5
$ws.Range("N107").Value = -12824584
  # Row 113: Copious Crystal Cannons
$ws.Range("H113").Value = 2018.8
$ws.Range("I113").Value = 1556.2941
$ws.Range("J113").Value = 3001.625
$ws.Range("K113").Value = 1556.2941
$ws.Range("L113").Value = 3001.625
$ws.Range("M113").Value = 613.7058999999999
$ws.Range("N113").Value = -7341.625

$ws = $wb.Worksheets.Item("LTW")
  # Row 2: Red in the Head
$ws.Range("H2").Value = 1062503.1
$ws.Range("I2").Value = 1100003.4
$ws.Range("J2").Value = 500000
$ws.Range("K2").Value = 1100003.4
$ws.Range("L2").Value = 500000
$ws.Range("M2").Value = -1099891.4
$ws.Range("N2").Value = -500224
  # Row 61: Spelling Me Softly
$ws.Range("H61").Value = 4592.933
$ws.Range("I61").Value = 1824.375
$ws.Range("K61").Value = 1824.375
$ws.Range("M61").Value = -1622.375
  # Row 113: Peace in Rest
$ws.Range("H113").Value = 4592.933
$ws.Range("I113").Value = 1824.375
$ws.Range("K113").Value = 1824.375
$ws.Range("M113").Value = 345.625
  # Row 132: Tenets of Tanning
$ws.Range("H132").Value = 710875.4399999999
$ws.Range("I132").Value = 862206.4
$ws.Range("K132").Value = 2586619.2
$ws.Range("M132").Value = -2584089.2
  # Row 136: Respect for Br'aax
$ws.Range("H136").Value = 1838.4445
$ws.Range("I136").Value = 1660.9231
$ws.Range("K136").Value = 4982.7693
$ws.Range("M136").Value = -2432.7693

$ws = $wb.Worksheets.Item("WVR")
  # Row 2: The Unmentionables
$ws.Range("H2").Value = 13612.571
$ws.Range("J2").Value = 17200
$ws.Range("L2").Value = 17200
$ws.Range("N2").Value = -17424
